$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates: D = Price (text), E = Volume(1h) percentage (text)
# D-column values are forced to Text format before assignment so Excel
# does not reinterpret dotted price strings (e.g. "322.03", "47.182.05")
# as numbers; the style is then reset to Normal so no stray numFmt/style
# survives on the cell (matching the source formatting).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.182.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.489.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.67%  "

$ws.Range("E7").Value = "  +0.80%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -1.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.878.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.493.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.847"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.101.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.24%  "

$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("E22").Value = "  +14.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "246.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.141"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.72%  "

$ws.Range("E32").Value = "  +0.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.53%  "

$ws.Range("E35").Value = "  +2.20%  "

$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("E37").Value = "  +1.88%  "

$ws.Range("E38").Value = "  +2.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.32%  "

$ws.Range("E40").Value = "  +0.36%  "

$ws.Range("E41").Value = "  -0.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "120.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.74%  "

$ws.Range("E44").Value = "  +0.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.992.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("E47").Value = "  -2.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.16%  "
